$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Core change: the "offset" row (row 2) used the error on the residual
# (C1 + $B$2) but should instead use the error on the mean (C1 - $B$2).
# Set the anchor cell (C2) on its own so it keeps a simple (non-shared)
# formula, then fill the rest of the row (D2:K2) as a single shared
# formula, matching the structure Excel originally used for this row.
$ws.Range("C2").Formula = '=C1-$B$2'
$ws.Range("D2:K2").Formula = '=D1-$B$2'

# --- Update the selected cell on Sheet1 (was D29, now B3) and make
# Sheet1 the active/selected sheet (it was Chart2 before) ---
$ws.Range("B3").Select()

$wb.Save()
